$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "Asset Ho Number" column values down one step to match the new
# inventory numbering scheme (PPAHOAP002->001, PPAHOAP003->002, PPAHOAP004->003)
$ws.Range("D17").Value = "PPAHOAP001"
$ws.Range("D18").Value = "PPAHOAP002"
$ws.Range("D19").Value = "PPAHOAP003"

# Header row re-wraps slightly with the refreshed content; match the new
# auto-computed row height.
$ws.Rows("16").RowHeight = 24.6
